# Mise à jour excel
# The boolean-ish columns (C:F = "est un garçon", "a les cheveux blonds",
# "a des lunettes", "a les cheveux roux") were stored as the text strings
# "true" / "false". Convert them to real numeric booleans (1/0) so the
# sheet no longer needs those two shared-string entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(1,1,0,0)
    3 = @(1,0,1,1)
    4 = @(0,1,1,0)
    5 = @(0,0,0,0)
    6 = @(1,0,1,0)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Range("C" + $row).Value = $values[0]
    $ws.Range("D" + $row).Value = $values[1]
    $ws.Range("E" + $row).Value = $values[2]
    $ws.Range("F" + $row).Value = $values[3]
}

# Move the active selection, matching the author's last cursor position.
$ws.Range("E23").Select()
